# Generate Report for Handoff
# Inserts two new "Ready for handoff" rows (749b72bf... and a5460bc7...)
# above the existing d4dc8520... row on each of the three sheets
# (Overview, zh-cn, de-de), pushing the d4dc8520... row down from row 3
# to row 5, and rewrites the hyperlinks that go with rows 3-5.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Drop existing hyperlinks on this sheet - Rows.Insert does not relocate
# a hyperlink's ref along with the cells it was anchored to, so we
# recreate all of them afterwards at their correct, final addresses.
$ws.Hyperlinks.Delete()

# Make room for the two new rows right above the current row 3
# (the d4dc8520... entry), which pushes it down to row 5.
$ws.Rows(3).Insert()
$ws.Rows(3).Insert()

$ws.Range("A3").Value = "749b72bf-41b4-4a05-bcd0-4eeabaf9f8dd.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-30-20 06:30:00"
$ws.Range("A3:D3").Style = $ws.Range("A5:D5").Style

$ws.Range("A4").Value = "a5460bc7-220b-4afd-a261-56e41741d2d2.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "2016-28-20 06:28:45"
$ws.Range("A4:D4").Style = $ws.Range("A5:D5").Style

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d50f05cd7d4681b034df6f614684e44eb2d83b36/e2e/77bde751-9604-4978-951d-bf3e7caef7fd.md", "", "", "77bde751-9604-4978-951d-bf3e7caef7fd.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3644bdd0d54a909153816295c4945464dd8293b/e2e/749b72bf-41b4-4a05-bcd0-4eeabaf9f8dd.md", "", "", "749b72bf-41b4-4a05-bcd0-4eeabaf9f8dd.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/a9bc4783fc2114c3b4ec75a025aba51d6ea1fdc/e2e/a5460bc7-220b-4afd-a261-56e41741d2d2.md", "", "", "a5460bc7-220b-4afd-a261-56e41741d2d2.md")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9bf7670fe5a96e1ee3bdbf48207781d40da0cece/e2e/d4dc8520-4bb8-4b24-87fb-2b204206540d.md", "", "", "d4dc8520-4bb8-4b24-87fb-2b204206540d.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Hyperlinks.Delete()

$ws.Rows(3).Insert()
$ws.Rows(3).Insert()
# Rows.Insert drags column F/G formatting down from row 2 - the source
# data never has those columns populated on the d4dc8520-style rows, so
# drop them again.
$ws.Range("F3:G4").Clear()

$ws.Range("A3").Value = "749b72bf-41b4-4a05-bcd0-4eeabaf9f8dd.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "749b72bf-41b4-4a05-bcd0-4eeabaf9f8dd.3644bdd0d54a909153816295c4945464dd8293ba.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-20 06:29:57"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"
$ws.Range("A3:I3").Style = $ws.Range("A5:I5").Style

$ws.Range("A4").Value = "a5460bc7-220b-4afd-a261-56e41741d2d2.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "a5460bc7-220b-4afd-a261-56e41741d2d2.a9bc4783fc2114c3b4ec75a025aba51d6ea1fdcd.zh-cn.xlf"
$ws.Range("E4").Value = "2016-03-20 06:29:57"
$ws.Range("H4").Value = "0001-01-01 00:00:00"
$ws.Range("I4").Value = "Include"
$ws.Range("A4:I4").Style = $ws.Range("A5:I5").Style

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d50f05cd7d4681b034df6f614684e44eb2d83b36/e2e/77bde751-9604-4978-951d-bf3e7caef7fd.md", "", "", "77bde751-9604-4978-951d-bf3e7caef7fd.md")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/d50f05cd7d4681b034df6f614684e44eb2d83b36/e2e/77bde751-9604-4978-951d-bf3e7caef7fd.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dba8594145e8f5fbfa7ce1a352c3b453a0e8bdb6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/77bde751-9604-4978-951d-bf3e7caef7fd.1febe2f7767a8891db474e9d95c5b92db5df0285.zh-cn.xlf", "", "", "77bde751-9604-4978-951d-bf3e7caef7fd.1febe2f7767a8891db474e9d95c5b92db5df0285.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/90bbb8dd1ab4f7495fc503e1fb4e3e8d2ff648dd/e2e/77bde751-9604-4978-951d-bf3e7caef7fd.md", "", "", "77bde751-9604-4978-951d-bf3e7caef7fd.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7f870959287fcac40545e89af5acd33d22498c19/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/77bde751-9604-4978-951d-bf3e7caef7fd.1febe2f7767a8891db474e9d95c5b92db5df0285.zh-cn.xlf", "", "", "77bde751-9604-4978-951d-bf3e7caef7fd.1febe2f7767a8891db474e9d95c5b92db5df0285.zh-cn.xlf")

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3644bdd0d54a909153816295c4945464dd8293b/e2e/749b72bf-41b4-4a05-bcd0-4eeabaf9f8dd.md", "", "", "749b72bf-41b4-4a05-bcd0-4eeabaf9f8dd.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/3644bdd0d54a909153816295c4945464dd8293b/e2e/749b72bf-41b4-4a05-bcd0-4eeabaf9f8dd.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3644bdd0d54a909153816295c4945464dd8293b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/749b72bf-41b4-4a05-bcd0-4eeabaf9f8dd.3644bdd0d54a909153816295c4945464dd8293ba.zh-cn.xlf", "", "", "749b72bf-41b4-4a05-bcd0-4eeabaf9f8dd.3644bdd0d54a909153816295c4945464dd8293ba.zh-cn.xlf")

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/a9bc4783fc2114c3b4ec75a025aba51d6ea1fdc/e2e/a5460bc7-220b-4afd-a261-56e41741d2d2.md", "", "", "a5460bc7-220b-4afd-a261-56e41741d2d2.md")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/a9bc4783fc2114c3b4ec75a025aba51d6ea1fdc/e2e/a5460bc7-220b-4afd-a261-56e41741d2d2.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a9bc4783fc2114c3b4ec75a025aba51d6ea1fdc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a5460bc7-220b-4afd-a261-56e41741d2d2.a9bc4783fc2114c3b4ec75a025aba51d6ea1fdcd.zh-cn.xlf", "", "", "a5460bc7-220b-4afd-a261-56e41741d2d2.a9bc4783fc2114c3b4ec75a025aba51d6ea1fdcd.zh-cn.xlf")

$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9bf7670fe5a96e1ee3bdbf48207781d40da0cece/e2e/d4dc8520-4bb8-4b24-87fb-2b204206540d.md", "", "", "d4dc8520-4bb8-4b24-87fb-2b204206540d.md")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/9bf7670fe5a96e1ee3bdbf48207781d40da0cece/e2e/d4dc8520-4bb8-4b24-87fb-2b204206540d.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/25038002da5bb3f8c1a9db11b27325e0a86f57fe/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d4dc8520-4bb8-4b24-87fb-2b204206540d.c3edfc3afc02d519f9f7643c04c820c1ccb51791.zh-cn.xlf", "", "", "d4dc8520-4bb8-4b24-87fb-2b204206540d.c3edfc3afc02d519f9f7643c04c820c1ccb51791.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Hyperlinks.Delete()

$ws.Rows(3).Insert()
$ws.Rows(3).Insert()
$ws.Range("F3:G4").Clear()

$ws.Range("A3").Value = "749b72bf-41b4-4a05-bcd0-4eeabaf9f8dd.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "749b72bf-41b4-4a05-bcd0-4eeabaf9f8dd.3644bdd0d54a909153816295c4945464dd8293ba.de-de.xlf"
$ws.Range("E3").Value = "2016-03-20 06:30:00"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"
$ws.Range("A3:I3").Style = $ws.Range("A5:I5").Style

$ws.Range("A4").Value = "a5460bc7-220b-4afd-a261-56e41741d2d2.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "a5460bc7-220b-4afd-a261-56e41741d2d2.a9bc4783fc2114c3b4ec75a025aba51d6ea1fdcd.de-de.xlf"
$ws.Range("E4").Value = "2016-03-20 06:30:00"
$ws.Range("H4").Value = "0001-01-01 00:00:00"
$ws.Range("I4").Value = "Include"
$ws.Range("A4:I4").Style = $ws.Range("A5:I5").Style

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d50f05cd7d4681b034df6f614684e44eb2d83b36/e2e/77bde751-9604-4978-951d-bf3e7caef7fd.md", "", "", "77bde751-9604-4978-951d-bf3e7caef7fd.md")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/d50f05cd7d4681b034df6f614684e44eb2d83b36/e2e/77bde751-9604-4978-951d-bf3e7caef7fd.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/20d73b9496eea8d892271f2bc2a0b76aa42aedcc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/77bde751-9604-4978-951d-bf3e7caef7fd.1febe2f7767a8891db474e9d95c5b92db5df0285.de-de.xlf", "", "", "77bde751-9604-4978-951d-bf3e7caef7fd.1febe2f7767a8891db474e9d95c5b92db5df0285.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8fdc0cb0b897bdca9b743420ccc2b6ea040ff5a6/e2e/77bde751-9604-4978-951d-bf3e7caef7fd.md", "", "", "77bde751-9604-4978-951d-bf3e7caef7fd.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1e5d078d21545d2f2326741c5e2c815258ff0241/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/77bde751-9604-4978-951d-bf3e7caef7fd.1febe2f7767a8891db474e9d95c5b92db5df0285.de-de.xlf", "", "", "77bde751-9604-4978-951d-bf3e7caef7fd.1febe2f7767a8891db474e9d95c5b92db5df0285.de-de.xlf")

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3644bdd0d54a909153816295c4945464dd8293b/e2e/749b72bf-41b4-4a05-bcd0-4eeabaf9f8dd.md", "", "", "749b72bf-41b4-4a05-bcd0-4eeabaf9f8dd.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/3644bdd0d54a909153816295c4945464dd8293b/e2e/749b72bf-41b4-4a05-bcd0-4eeabaf9f8dd.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3644bdd0d54a909153816295c4945464dd8293b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/749b72bf-41b4-4a05-bcd0-4eeabaf9f8dd.3644bdd0d54a909153816295c4945464dd8293ba.de-de.xlf", "", "", "749b72bf-41b4-4a05-bcd0-4eeabaf9f8dd.3644bdd0d54a909153816295c4945464dd8293ba.de-de.xlf")

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/a9bc4783fc2114c3b4ec75a025aba51d6ea1fdc/e2e/a5460bc7-220b-4afd-a261-56e41741d2d2.md", "", "", "a5460bc7-220b-4afd-a261-56e41741d2d2.md")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/a9bc4783fc2114c3b4ec75a025aba51d6ea1fdc/e2e/a5460bc7-220b-4afd-a261-56e41741d2d2.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a9bc4783fc2114c3b4ec75a025aba51d6ea1fdc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a5460bc7-220b-4afd-a261-56e41741d2d2.a9bc4783fc2114c3b4ec75a025aba51d6ea1fdcd.de-de.xlf", "", "", "a5460bc7-220b-4afd-a261-56e41741d2d2.a9bc4783fc2114c3b4ec75a025aba51d6ea1fdcd.de-de.xlf")

$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9bf7670fe5a96e1ee3bdbf48207781d40da0cece/e2e/d4dc8520-4bb8-4b24-87fb-2b204206540d.md", "", "", "d4dc8520-4bb8-4b24-87fb-2b204206540d.md")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/9bf7670fe5a96e1ee3bdbf48207781d40da0cece/e2e/d4dc8520-4bb8-4b24-87fb-2b204206540d.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8ddd0d4175dd7ddf2b648d8caae689eee7cbaabf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d4dc8520-4bb8-4b24-87fb-2b204206540d.c3edfc3afc02d519f9f7643c04c820c1ccb51791.de-de.xlf", "", "", "d4dc8520-4bb8-4b24-87fb-2b204206540d.c3edfc3afc02d519f9f7643c04c820c1ccb51791.de-de.xlf")
